$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): add P1=14, Q1=15, matching the style of O1 ---
$ws.Cells.Item(1, 15).Copy()
$ws.Cells.Item(1, 16).PasteSpecial(-4122)
$ws.Cells.Item(1, 16).Value = 14

$ws.Cells.Item(1, 15).Copy()
$ws.Cells.Item(1, 17).PasteSpecial(-4122)
$ws.Cells.Item(1, 17).Value = 15

$excel.CutCopyMode = 0

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Column I: 1 -> 2
    $ws.Cells.Item($r, 9).Value = 2
    # Column K: 2 -> 1
    $ws.Cells.Item($r, 11).Value = 1
    # Column M: 1 -> 2
    $ws.Cells.Item($r, 13).Value = 2
    # Column O: 2 -> 1
    $ws.Cells.Item($r, 15).Value = 1
    # New columns P and Q, value 2 each
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
